$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remembered Excel window position/size (best effort; some hosts do not
# persist these window-chrome values back into workbookView) ---
$w1 = $wb.Windows.Item(1)
$w1.Left = 2140
$w1.Top = 3240
$w1.Width = 13820
$w1.Height = 16860

# --- New data rows (5-12), appended after the existing 4 data rows ---
# Row 5 - "one"
$ws.Range("A5").Value = "one"
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 18

# Row 6 - "three"
$ws.Range("A6").Value = "three"
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 39
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 9

# Row 7 - "four"
$ws.Range("A7").Value = "four"
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9

# Row 8 - "five"
$ws.Range("A8").Value = "five"
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 39
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 17
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 10

# Row 10 - "v3_c"
$ws.Range("A10").Value = "v3_c"
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 40
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11

# Row 11 - "v3_d"
$ws.Range("A11").Value = "v3_d"
$ws.Range("B11").Value = 17
$ws.Range("C11").Value = 38
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 18
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 10

# Row 12 - "v3_e"
$ws.Range("A12").Value = "v3_e"
$ws.Range("B12").Value = 18
$ws.Range("C12").Value = 40
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 10

# Row 9 - "v3_a" (added after v3_c/d/e, matching shared-string allocation order)
$ws.Range("A9").Value = "v3_a"
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 40
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 18

# --- New "Reason"/"Magnitude" annotation column (H) ---
$ws.Range("H1").Value = "Reason"
$ws.Range("H5").Value = "Magnitude (x)"
$ws.Range("H3").Value = "Magnitude (x,y)"
$ws.Range("H7").Value = "Magnitude (x,y)"
$ws.Range("H8").Value = "Magnitude (x,y)"
$ws.Range("H12").Value = "Magnitude (x,y)"

# Column H width
$ws.Columns.Item(8).ColumnWidth = 14

# --- Highlight subject-label cells in column A ---
# Establish fill color creation order: green, then yellow, then red
$ws.Range("A2").Interior.Color = 5287936   # green FF00B050
$ws.Range("A4").Interior.Color = 65535     # yellow FFFFFF00
$ws.Range("A3").Interior.Color = 255       # red FFFF0000

$ws.Range("A6").Interior.Color = 5287936   # green
$ws.Range("A9").Interior.Color = 5287936   # green
$ws.Range("A10").Interior.Color = 5287936  # green

$ws.Range("A11").Interior.Color = 65535    # yellow

$ws.Range("A5").Interior.Color = 255       # red
$ws.Range("A7").Interior.Color = 255       # red
$ws.Range("A8").Interior.Color = 255       # red
$ws.Range("A12").Interior.Color = 255      # red

# --- Selection / active cell ---
$ws.Range("F4").Select()
